$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "DVD - Other-Right (Divers droite, DVD)"
$ws.Range("C1").Value = "DVG - Other-Left (Divers gauche, DVG)"
$ws.Range("D1").Value = "PC - Communist Party (Parti communiste, PC)"
$ws.Range("E1").Value = "PS - Socialist Party (Parti socialiste, PS)"
$ws.Range("F1").Value = "RPR - Rally for the Republic (Rassemblement pour la République, RPR)"
$ws.Range("G1").Value = "UDF - Union for French Democracy (Union pour la Démocratie française, UDF)"
$ws.Range("H1").Value = "FN - National Front (Front national, FN)"
$ws.Range("I1").Value = "V - Greens (Les Verts, V)"
$ws.Range("J1").Value = "PR - Republican Alliance (Pôle républicain, PR)"
$ws.Range("K1").Value = "UMP, LR - Union for a Popular Movement (Union pour un mouvement populaire, UMP, LR), known until  as Union for a Presidential Majority (UMP, Union pour la majorité présidentielle)"
$ws.Range("L1").Value = "DREG - Other-Regionalists (Régionaliste, DREG)"
$ws.Range("M1").Value = "MPF - Movement for France (Mouvement pour la France, MPF)"
$ws.Range("N1").Value = "NC - New Centre (Nouveau centre, NC)"
$ws.Range("O1").Value = "Other - - (-, Other)"
$ws.Range("P1").Value = "RdG, PRG - Left Radicals (Parti des Radicaux de gauche, RdG, PRG)"
$ws.Range("Q1").Value = "AC - Centrist Alliance (Alliance centriste , AC)"
$ws.Range("R1").Value = "EXD - Extreme Right (Extrême droite, EXD)"
$ws.Range("S1").Value = "FdG - Left Front/Alliance of the Overseas (Front de gauche, FdG)"
$ws.Range("T1").Value = "Radical - Radical Party (Parti radical, Radical)"
$ws.Range("U1").Value = "D - Diverse (Divers, D)"
$ws.Range("V1").Value = "DLF - France Arise (Debout la France, DLF)"
$ws.Range("W1").Value = "DVE - Other-Ecologists (Autres écologistes, DVE)"
$ws.Range("X1").Value = "FI - Unsubmissive France (La France Insoumise, FI)"
$ws.Range("Y1").Value = "LRM - Forward (La Republique en marche, LRM)"
$ws.Range("Z1").Value = "MoDem - Democratic Movement (Mouvement démocrate, MoDem)"
